$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-font the existing header + data rows (1:36) from the Khmer-specific
#    fonts to Arial (keeping the same point size / colour roles) and widen
#    the key/value/translation columns so the longer strings fit.
# ---------------------------------------------------------------------------
$ws.Columns("A:C").Font.Name = "Arial"

$ws.Columns("A").ColumnWidth = 41.67
$ws.Columns("B").ColumnWidth = 86.34
$ws.Columns("C").ColumnWidth = 97.67

$ws.Rows("1:36").RowHeight = 29.25

# ---------------------------------------------------------------------------
# 2. Append the new translation keys introduced by this commit (save /
#    delete / user / time / month / to / from) as rows 37-43.
# ---------------------------------------------------------------------------
$newRows = @(
  @("save", "Save", "រក្សាទុក"),
  @("delete", "Delete", "លុប"),
  @("user", "User", "អ្នកប្រើប្រាស់"),
  @("time", "Time", "ម៉ោង"),
  @("month", "Month", "ខែ"),
  @("to", "to", "ទៅ"),
  @("from", "from", "ពី")
)

$r = 37
foreach ($row in $newRows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Rows.Item($r).RowHeight = 29.25
  $r++
}

$ws.Range("A37:C43").Font.Name = "Arial"

# ---------------------------------------------------------------------------
# 3. Restore the view: keep the header frozen, scroll near the newly added
#    rows and leave the selection on the last value cell, same as the author
#    left it after typing the new rows in.
# ---------------------------------------------------------------------------
$ws.Range("B35").Select()
